# Applies the forecasts_table_SP.xlsx edit:
#   - extends both "cases" and "deaths" sheets with two new forecast-origin
#     columns AC ("2020-05-08") and AD ("2020-05-09"), and two new observed
#     rows 41 ("2020-05-22") and 42 ("2020-05-23")
#   - backfills previously-empty B27/B28 "Observed" cells
#   - fills in the new AC/AD forecast values for rows 27-42
#
# NOTE on text cells: the source date-like labels ("2020-05-08", "2020-05-22",
# ...) must land in the sheet as TEXT (shared strings), matching how the
# workbook already stores every other date-label column/row header - not as
# Excel auto-converted date serial numbers. Plain `.Value = "2020-05-08"`
# assignment is auto-parsed by COM into a date serial, so:
#   * for labels that already exist elsewhere in the sheet (the new headers
#     AC1/AD1 duplicate the existing A27/A28 row-label text), we `.Copy()`
#     the existing text cell onto the new cell - this clones the shared
#     string + (lack of) style exactly, no auto date parsing involved.
#   * for the genuinely new labels (new row headers A41/A42, which need two
#     brand new shared strings) we type them with a leading apostrophe, the
#     standard Excel "force text" prefix, then propagate via `.Copy()` to
#     the other sheet so the trick only runs once per label.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("cases")
$ws2 = $wb.Worksheets.Item("deaths")

# --- new header cells AC1 / AD1 on both sheets -----------------------------
# "2020-05-08" already lives at A27 (t="s") and "2020-05-09" at A28; copying
# those cells onto AC1/AD1 reuses the existing shared-string entries.
$ws1.Cells.Item(27, 1).Copy($ws1.Cells.Item(1, 29))
$ws1.Cells.Item(28, 1).Copy($ws1.Cells.Item(1, 30))
$ws2.Cells.Item(27, 1).Copy($ws2.Cells.Item(1, 29))
$ws2.Cells.Item(28, 1).Copy($ws2.Cells.Item(1, 30))

# --- new row-label cells A41 / A42 on both sheets --------------------------
# These two labels ("2020-05-22", "2020-05-23") are brand new strings, so
# type them once (forcing text with a leading apostrophe) and copy across.
# The apostrophe prefix leaves a "quote prefix" flag on the cell's style;
# resetting to the "Normal" cell style afterwards drops that flag again so
# the cell ends up with the same (default) style as every other label cell.
$ws1.Cells.Item(41, 1).Value = "'2020-05-22"
$ws1.Cells.Item(41, 1).Style = "Normal"
$ws1.Cells.Item(42, 1).Value = "'2020-05-23"
$ws1.Cells.Item(42, 1).Style = "Normal"
$ws1.Cells.Item(41, 1).Copy($ws2.Cells.Item(41, 1))
$ws1.Cells.Item(42, 1).Copy($ws2.Cells.Item(42, 1))

# --- numeric cell values: [row, col, value] ---------------------------------
# cases: B27/B28 backfill + AC/AD forecast columns for rows 27-42
$casesValues = @(
    @(27,2,41830),
    @(28,2,44411),
    @(28,29,43290),
    @(29,29,44416),
    @(29,30,45560),
    @(30,29,45866),
    @(30,30,46911),
    @(31,29,47085),
    @(31,30,48155),
    @(32,29,48019),
    @(32,30,49204),
    @(33,29,49150),
    @(33,30,50318),
    @(34,29,49914),
    @(34,30,51224),
    @(35,29,51096),
    @(35,30,52414),
    @(36,29,51807),
    @(36,30,53194),
    @(37,29,52708),
    @(37,30,54095),
    @(38,29,53546),
    @(38,30,54954),
    @(39,29,54480),
    @(39,30,55870),
    @(40,29,55565),
    @(40,30,56886),
    @(41,29,56500),
    @(41,30,57843),
    @(42,30,58745)
)
foreach ($item in $casesValues) {
    $ws1.Cells.Item($item[0], $item[1]).Value = $item[2]
}

# deaths: same layout, different figures
$deathsValues = @(
    @(27,2,3416),
    @(28,2,3608),
    @(28,29,3504),
    @(29,29,3618),
    @(29,30,3684),
    @(30,29,3722),
    @(30,30,3781),
    @(31,29,3819),
    @(31,30,3877),
    @(32,29,3906),
    @(32,30,3963),
    @(33,29,3985),
    @(33,30,4044),
    @(34,29,4058),
    @(34,30,4114),
    @(35,29,4126),
    @(35,30,4183),
    @(36,29,4199),
    @(36,30,4251),
    @(37,29,4262),
    @(37,30,4311),
    @(38,29,4327),
    @(38,30,4371),
    @(39,29,4395),
    @(39,30,4431),
    @(40,29,4460),
    @(40,30,4495),
    @(41,29,4520),
    @(41,30,4556),
    @(42,30,4612)
)
foreach ($item in $deathsValues) {
    $ws2.Cells.Item($item[0], $item[1]).Value = $item[2]
}
